$d = $word.ActiveDocument

# Locate the target paragraph ("Create two separate Scenes in Unity and
# build them to your Android mobile device.") by scanning the paragraph
# collection instead of hard-coding character offsets, so the script is
# resilient to any earlier content changes.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Create two separate Scenes in Unity and build them to your Android mobile device.*") {
        $target = $p
        break
    }
}

$paraStart = $target.Range.Start
$paraEnd = $target.Range.End

# Bold "Create" (offsets are relative to the paragraph start).
$r1 = $d.Range($paraStart + 0, $paraStart + 6)
$r1.Bold = 1

# Bold "and build them".
$r2 = $d.Range($paraStart + 36, $paraStart + 50)
$r2.Bold = 1
